# Updated cryptos list on Wed Aug 30 02:04:20 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for every coin row,
# and swaps the Maker / TrustWalletToken rows (40 and 42) to reflect the
# new ranking order, including their Coin name / Link / Price / Volume.
#
# Note: several price strings look numeric (e.g. "0.8600", "0.06630").
# Excel would otherwise auto-convert such text into a Number and silently
# drop the trailing zeros, so NumberFormat is forced to Text ("@") for
# those specific cells right before the value is written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.647.02'
$ws.Range("E2").Value = '  +5.50%  '
$ws.Range("D3").Value = '1.724.46'
$ws.Range("E3").Value = '  +3.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.35'
$ws.Range("E5").Value = '  +3.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5382'
$ws.Range("E6").Value = '  +2.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2707'
$ws.Range("E8").Value = '  +1.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06630'
$ws.Range("E9").Value = '  +3.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.76'
$ws.Range("E10").Value = '  +4.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07761'
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.657'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("D13").Value = '1.717.04'
$ws.Range("E13").Value = '  +4.28%  '
$ws.Range("D14").Value = '1.962.09'
$ws.Range("E14").Value = '  +3.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5899'
$ws.Range("E15").Value = '  +4.40%  '
$ws.Range("D16").Value = '0.0₅8312'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.16'
$ws.Range("E17").Value = '  +3.62%  '
$ws.Range("D18").Value = '27.639.37'
$ws.Range("E18").Value = '  +5.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '225.29'
$ws.Range("E19").Value = '  +16.58%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.763'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.75'
$ws.Range("E22").Value = '  +1.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.136'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.77'
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.702'
$ws.Range("E26").Value = '  +11.91%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1237'
$ws.Range("E27").Value = '  +2.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.434'
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.85'
$ws.Range("E29").Value = '  +4.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05591'
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.306'
$ws.Range("E31").Value = '  +2.21%  '
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.477'
$ws.Range("E33").Value = '  +2.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.665'
$ws.Range("E34").Value = '  +5.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9667'
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.446'
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5949'
$ws.Range("E38").Value = '  +3.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01657'
$ws.Range("E39").Value = '  +3.46%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.070.95'
$ws.Range("E40").Value = '  +2.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.862'
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8600'
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.37'
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").Value = '1.867.64'
$ws.Range("E45").Value = '  +3.71%  '
$ws.Range("E46").Value = '  +8.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '59.34'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.225'
$ws.Range("E48").Value = '  +1.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4429'
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9987'
$ws.Range("E50").Value = '  -0.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05288'
$ws.Range("E51").Value = '  +0.74%  '
